{"js": "// The document was originally authored with one run per word (and one run\n// per separating space). This edit consolidates each of those runs into a\n// single run holding the full paragraph text, without changing any of the\n// visible text itself. It touches the three opening paragraphs: the Title,\n// the Author line, and the Abstract body (the AbstractTitle \"Summary\"\n// paragraph already had a single run and needs no change).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nconst replacements = {\n  \"Title\": \"Questions: Introduction to differentiation and the derivative\",\n  \"Author\": \"Sara Delgado Garcia\",\n  \"Abstract\": \"A selection of questions for the study guide on introduction to differentiation and the derivative.\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const style = para.style;\n  if (Object.prototype.hasOwnProperty.call(replacements, style)) {\n    para.insertText(replacements[style], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document was originally authored with one run per word (plus one run\n# per separating space). This edit consolidates each of those runs into a\n# single run holding the full paragraph text, without changing any of the\n# visible text itself. It touches the three opening paragraphs: the Title,\n# the Author line, and the Abstract body (the AbstractTitle \"Summary\"\n# paragraph already had a single run and needs no change).\n#\n# Note: assigning Range.Text on a range that spans several runs only\n# rewrites the first run in place here, leaving the remaining runs (and\n# the paragraph mark) untouched. So after setting the new text we compute\n# the leftover tail range (from right after the freshly written text up to\n# just before the paragraph mark) and delete it, which leaves the\n# paragraph holding a single run with exactly the desired text.\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($para, [string]$newText) {\n    $paraStart = $para.Range.Start\n    $para.Range.Text = $newText\n    $afterEnd = $para.Range.End\n    $tail = $d.Range($paraStart + $newText.Length, $afterEnd - 1)\n    if ($tail.Start -lt $tail.End) {\n        $tail.Delete()\n    }\n}\n\n$replacements = @{\n    \"Title\"    = \"Questions: Introduction to differentiation and the derivative\"\n    \"Author\"   = \"Sara Delgado Garcia\"\n    \"Abstract\" = \"A selection of questions for the study guide on introduction to differentiation and the derivative.\"\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $styleName = $para.Style.NameLocal\n    if ($replacements.ContainsKey($styleName)) {\n        Set-ParagraphText $para $replacements[$styleName]\n    }\n}\n"}
